# Helper: write a value as TEXT (not auto-converted to a number) while
# keeping the cell's existing style/format index (e.g. no border/bold
# carried over from the "@" number format we briefly apply).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计", by duplicating
#    the existing "2022-Q3" sheet (identical header layout/styling)
#    and then trimming it down to a single data row with new values.
# ------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Copy($null, $zongji)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop all the data rows that came from 2022-Q3 except the first one,
# which we will overwrite with the 2022-Q4 fund row below.
$q4.Rows("3:14").Delete() | Out-Null

$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "010434"
$q4.Range("C2").Value = "红土创新医疗保健股票"
Set-TextValue $q4.Range("D2") "3.12"
Set-TextValue $q4.Range("E2") "94.35"
Set-TextValue $q4.Range("F2") "7.05"
Set-TextValue $q4.Range("G2") "0.2200"
$q4.Range("H2").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add the 2022-Q4 row at the top
#    of the data, push the other rows down one slot, and renumber the
#    index column (A) sequentially 0..4.
# ------------------------------------------------------------------
$zongji.Range("A5").Copy($zongji.Range("A6"))

$zongji.Range("A6").Value = 4
$zongji.Range("B6").Value = "2021-Q4"
$zongji.Range("C6").Value = 2
$zongji.Range("D6").Value = 0.3

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2022-Q1"
$zongji.Range("C5").Value = 7
$zongji.Range("D5").Value = 0.63

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2022-Q2"
$zongji.Range("C4").Value = 12
$zongji.Range("D4").Value = 0.58

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q3"
$zongji.Range("C3").Value = 13
$zongji.Range("D3").Value = 0.49

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 1
$zongji.Range("D2").Value = 0.22

# Restore "总计" as the active sheet (creating/copying a sheet makes it
# active, same as it would in the Excel UI) to match the original
# workbook's selection state.
$zongji.Activate()
